# Column B ("交易日期" / Transaction Date) stores dates as plain
# YYYYMMDD integers (e.g. 20190202). This script rewrites every data row
# (B2:B304) with the equivalent Excel date-serial number (e.g. 43498) and
# applies a "YYYY-MM-DD" date display format to the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "row:serial" pairs for every data row, in row order.
$data = @("2:43498","3:43498","4:43498","5:43508","6:43509","7:43509","8:43511","9:43526","10:43526","11:43535","12:43535","13:43535","14:43535","15:43535","16:43535","17:43535","18:43535","19:43535","20:43536","21:43536","22:43538","23:43538","24:43538","25:43542","26:43543","27:43544","28:43545","29:43545","30:43545","31:43545","32:43545","33:43545","34:43545","35:43545","36:43545","37:43545","38:43557","39:43557","40:43557","41:43563","42:43563","43:43563","44:43565","45:43565","46:43565","47:43565","48:43565","49:43565","50:43565","51:43565","52:43565","53:43565","54:43566","55:43587","56:43587","57:43587","58:43593","59:43595","60:43595","61:43595","62:43595","63:43595","64:43595","65:43595","66:43595","67:43595","68:43595","69:43595","70:43595","71:43595","72:43595","73:43595","74:43595","75:43595","76:43595","77:43595","78:43595","79:43598","80:43605","81:43605","82:43608","83:43608","84:43608","85:43608","86:43609","87:43609","88:43614","89:43614","90:43614","91:43614","92:43614","93:43618","94:43618","95:43618","96:43620","97:43621","98:43622","99:43626","100:43626","101:43626","102:43626","103:43626","104:43626","105:43626","106:43626","107:43626","108:43626","109:43626","110:43626","111:43626","112:43626","113:43626","114:43626","115:43626","116:43627","117:43628","118:43636","119:43637","120:43641","121:43643","122:43644","123:43644","124:43648","125:43648","126:43648","127:43649","128:43654","129:43654","130:43654","131:43654","132:43656","133:43656","134:43656","135:43656","136:43656","137:43656","138:43656","139:43656","140:43656","141:43656","142:43656","143:43656","144:43657","145:43658","146:43658","147:43658","148:43663","149:43671","150:43677","151:43677","152:43678","153:43678","154:43679","155:43679","156:43679","157:43683","158:43683","159:43686","160:43686","161:43686","162:43686","163:43686","164:43686","165:43686","166:43686","167:43686","168:43686","169:43689","170:43690","171:43691","172:43696","173:43697","174:43707","175:43707","176:43707","177:43710","178:43710","179:43710","180:43717","181:43717","182:43718","183:43718","184:43718","185:43718","186:43718","187:43718","188:43718","189:43718","190:43718","191:43718","192:43718","193:43718","194:43718","195:43719","196:43719","197:43719","198:43729","199:43737","200:43737","201:43737","202:43740","203:43740","204:43740","205:43748","206:43748","207:43748","208:43748","209:43748","210:43748","211:43748","212:43748","213:43748","214:43749","215:43750","216:43756","217:43756","218:43756","219:43756","220:43762","221:43768","222:43768","223:43768","224:43768","225:43768","226:43770","227:43771","228:43771","229:43771","230:43772","231:43775","232:43776","233:43777","234:43777","235:43777","236:43777","237:43780","238:43780","239:43780","240:43781","241:43787","242:43787","243:43789","244:43801","245:43801","246:43801","247:43801","248:43801","249:43802","250:43809","251:43809","252:43809","253:43809","254:43809","255:43809","256:43809","257:43809","258:43809","259:43809","260:43809","261:43810","262:43811","263:43820","264:43826","265:43826","266:43826","267:43832","268:43832","269:43832","270:43836","271:43840","272:43840","273:43840","274:43840","275:43840","276:43840","277:43840","278:43840","279:43840","280:43840","281:43843","282:43843","283:43843","284:43849","285:43849","286:43849","287:43863","288:43863","289:43863","290:43871","291:43871","292:43871","293:43871","294:43872","295:43879","296:43879","297:43888","298:43889","299:43889","300:43889","301:43889","302:43892","303:43892","304:43893")

foreach ($entry in $data) {
    $parts  = $entry.Split(":")
    $row    = [int]$parts[0]
    $serial = [double]$parts[1]
    $ws.Range("B$row").Value = $serial
}

# Register the "YYYY-MM-DD" number format and give B2 a fresh cell style.
$ws.Range("B2").NumberFormat = "yyyy-mm-dd"
$ws.Range("B2").NumberFormat = "YYYY-MM-DD"

# Apply the same date format to the rest of the column.
foreach ($entry in $data) {
    $parts = $entry.Split(":")
    $row   = [int]$parts[0]
    if ($row -ne 2) {
        $ws.Range("B$row").NumberFormat = "YYYY-MM-DD"
    }
}

Write-Host "Converted $($data.Count) dates in column B to date-serial values with YYYY-MM-DD formatting"
